$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the stray empty paragraph that used to follow the
#    "М.А.Маренный" signature line. (Do this before touching the
#    table so paragraph indices stay sane.)
# ------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "`r") {
        $prev = $paras.Item($i - 1)
        if ($prev.Range.Text -like "*Маренный*") {
            $p.Range.Delete()
            break
        }
    }
}

# ------------------------------------------------------------------
# 2. The info table: shrink every paragraph's font from 12pt (sz 24)
#    to 11pt (sz 22) - both the run and the paragraph-mark formatting.
# ------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$tbl.Range.Font.Size = 11
$tbl.Range.Font.SizeBi = 11

# ------------------------------------------------------------------
# 3. Address cell: drop ALL CAPS character formatting, but bake the
#    upper-cased text directly into the run, and mark the edit spot
#    with the "_GoBack" bookmark (as Word does for the last edit).
# ------------------------------------------------------------------
$addrCell = $tbl.Cell(2, 2)
$addrRange = $addrCell.Range
$addrRange.Font.AllCaps = 0
$addrRange.Text = "117513, Г.МОСКВА, ВН.ТЕР.Г. МУНИЦИПАЛЬНЫЙ ОКРУГ ТЕПЛЫЙ СТАН, УЛ.ОСТРОВИТЯНОВА, Д.6, ПОМЕЩ. 3/П"

$findRange = $d.Content
$findRange.Find.Execute("117513, Г.МОСКВА, ВН.ТЕР.Г. МУНИЦИПАЛЬНЫЙ ОКРУГ ТЕПЛЫЙ СТАН, УЛ.ОСТРОВИТЯНОВА, Д.6, ПОМЕЩ. 3/П", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $findRange)

# ------------------------------------------------------------------
# 4. Remove the previous "_GoBack" bookmark location (it used to sit
#    in the first-page footer, right before the signature text box).
# ------------------------------------------------------------------
$sec = $d.Sections(1)
$ftr = $sec.Footers(2)
$ftrRange = $ftr.Range
$hasOldBookmark = $ftrRange.Bookmarks.Exists("_GoBack")
if ($hasOldBookmark) {
    $oldBm = $ftrRange.Bookmarks("_GoBack")
    $oldBm.Range.Delete()
}
